# Hours Worked.xlsx - add week commencing 21/05/2018 (lap 3), with an
# "X" marker in the Friday/Saturday/Sunday columns, plus a couple of
# missed hours entries for the week before, and move the Total row down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Week of 30/04/2018 (row 16): extra hours recorded for Saturday/Sunday
$ws.Range("G16").Value = 5
$ws.Range("H16").Value = 3

# --- Week of 07/05/2018 (row 17): Tuesday/Wednesday entries removed (blank instead of 0)
$ws.Range("C17").ClearContents()
$ws.Range("D17").ClearContents()

# --- Week of 14/05/2018 (row 18): Tuesday/Wednesday entries removed (blank instead of 0)
$ws.Range("C18").ClearContents()
$ws.Range("D18").ClearContents()
$ws.Range("I18").Formula = "=SUM(B18:H18)"

# --- Insert a new row for week of 21/05/2018, pushing the Total row down
$ws.Rows("19").Insert()

$ws.Range("A19").Value = 43241
$ws.Range("F19").Value = "X"
$ws.Range("G19").Value = "X"
$ws.Range("H19").Value = "X"
$ws.Range("I19").Formula = "=SUM(B19:E19)"

# --- Total row, now on row 20, sums through the newly added row 19
$ws.Range("I20").Formula = "=SUM(I2:I19)"

# --- Restore/update the selection shown when the workbook was last saved
$ws.Range("C17:D19").Select()
